$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.395.05"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.876.10"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7134"
$ws.Range("E5").Value = "  -1.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.80"
$ws.Range("E6").Value = "  +0.46%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3114"
$ws.Range("E8").Value = "  +0.88%  "
$ws.Range("E9").Value = "  -2.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.26"
$ws.Range("E10").Value = "  +0.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08384"
$ws.Range("E11").Value = "  +1.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.902.73"
$ws.Range("E12").Value = "  +1.99%  "
$ws.Range("E13").Value = "  +0.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7163"
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.75"
$ws.Range("E15").Value = "  +1.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.406.39"
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008302"
$ws.Range("E17").Value = "  +6.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.973"
$ws.Range("E18").Value = "  +2.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.85"
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.136.39"
$ws.Range("E20").Value = "  +1.05%  "
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9996"
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.911"
$ws.Range("E23").Value = "  -0.78%  "
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1618"
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.88"
$ws.Range("E26").Value = "  +1.07%  "
$ws.Range("E27").Value = "  +0.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.60"
$ws.Range("E28").Value = "  +1.99%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.508"
$ws.Range("E29").Value = "  +0.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.414"
$ws.Range("E30").Value = "  +1.00%  "
$ws.Range("E31").Value = "  -4.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.329"
$ws.Range("E32").Value = "  +5.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05227"
$ws.Range("E33").Value = "  +0.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.928"
$ws.Range("E34").Value = "  -0.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7754"
$ws.Range("E35").Value = "  +6.83%  "
$ws.Range("E36").Value = "  -0.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.680"
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01865"
$ws.Range("E38").Value = "  +0.48%  "
$ws.Range("E39").Value = "  +0.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.163.54"
$ws.Range("E40").Value = "  -0.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.401"
$ws.Range("E41").Value = "  +4.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "73.60"
$ws.Range("E42").Value = "  +1.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8912"
$ws.Range("E43").Value = "  -1.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "104.57"
$ws.Range("E44").Value = "  +2.79%  "
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.029.93"
$ws.Range("E46").Value = "  +0.92%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5203"
$ws.Range("E47").Value = "  -1.48%  "
$ws.Range("E48").Value = "  +0.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.403"
$ws.Range("E49").Value = "  +1.56%  "
$ws.Range("E50").Value = "  +2.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4311"
$ws.Range("E51").Value = "  +0.84%  "
